# "added two more problems"
# Append a new row (row 14) to the problems table with a new section header,
# problem name and solution/approach notes, following the same layout and
# formatting conventions used by the existing "Fast & Slow Pointers" rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Values for the new row -------------------------------------------------
$ws.Range("A14").Value = "14. Fast & Slow Pointers"
$ws.Range("B14").Value = "Remove duplicates"
$ws.Range("C14").Value = "Travers through the linked list. If curr.next "

# --- Formatting --------------------------------------------------------------
# Column B (problem name) in this section uses the same highlighted fill as
# the rest of the "Fast & Slow Pointers" block (e.g. B13). Column C (notes)
# uses the plain "vertical top" style used by several rows (e.g. C8).
$ws.Range("B13").Copy() | Out-Null
$ws.Range("B14").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats, `
    [Microsoft.Office.Interop.Excel.XlPasteSpecialOperation]::xlPasteSpecialOperationNone, `
    $false, $false) | Out-Null

$ws.Range("C8").Copy() | Out-Null
$ws.Range("C14").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats, `
    [Microsoft.Office.Interop.Excel.XlPasteSpecialOperation]::xlPasteSpecialOperationNone, `
    $false, $false) | Out-Null

$excel.CutCopyMode = 0

# Row height matches the new row's shorter content.
$ws.Rows.Item(14).RowHeight = 22.5

# Leave the selection where it was left after entering the last value.
$ws.Range("C14").Select() | Out-Null
